# Add "COG Dataset 6" block to the correlation pipeline.
# For each of the 4 worksheets, insert three new rows right before the
# existing "fMRI Dataset" block (old rows 10-12) and populate them with
# the new COG Dataset 6 correlation figures. The old rows 10-12 shift
# down to rows 13-15 automatically as part of the row insert.

$wb = $excel.ActiveWorkbook

$newLabels = @(
    "Time (in sec.) - COG Dataset 6",
    "Correctness (% of correct answers given to verification Qs) - COG Dataset 6",
    "Rating (True or False) - COG Dataset 6"
)

# Per-sheet (by tab order) data for the new rows 10, 11, 12 -> columns B..I
$sheetData = @{
    1 = @(
        @(50, 39, 192, 50, -0.2516122567491849,  0.02940004335865699, -0.311094163299463,   0.02787696039183273),
        @(50, 39, 192, 50, -0.02439922062337898,  0.8337884478927275, -0.02888208251072132,  0.8421820698041393),
        @(50, 39, 192, 50, -0.02584053907657639,  0.8319012317886401, -0.02491868767248899,  0.8636170341288709)
    )
    2 = @(
        @(50, 28, 83, 50, -0.2516122567491849,   0.02940004335865699, -0.3118320089283974,   0.02748882596687927),
        @(50, 28, 83, 50, -0.02439922062337898,   0.8337884478927275, -0.02805160729475809,  0.8466646412386398),
        @(50, 28, 83, 50, -0.01788960397609135,   0.8831700141519032, -0.02174552145202223,  0.8808493755591094)
    )
    3 = @(
        @(50, 14, 86, 50, -0.2516122567491849,   0.02940004335865699, -0.311094163299463,    0.02787696039183273),
        @(50, 14, 86, 50, -0.01689176812387775,   0.8844917889522957, -0.02140780556705223,  0.8826867552381472),
        @(50, 14, 86, 50, -0.02584053907657639,   0.8319012317886401, -0.02445204558124388,  0.8661474618649031)
    )
    4 = @(
        @(50, 22, 23, 50, -0.2513580286171832,   0.03054279092665989, -0.3106508608887624,   0.02811234624135277),
        @(50, 22, 23, 50, -0.02644542653481615,   0.8211283092369018, -0.02889232075369115,  0.8421268384730887),
        @(50, 22, 23, 50, -0.008002176888210664,  0.9479027969797296, -0.008309173651488432, 0.9543304334608099)
    )
}

for ($s = 1; $s -le 4; $s++) {
    $ws = $wb.Worksheets.Item($s)

    # Insert three fresh rows above the current row 10 ("fMRI Dataset" block),
    # pushing that block down to rows 13-15. (Insert row 10 three times,
    # rather than a "10:12" range insert, for reliable shifting.)
    $ws.Rows.Item(10).Insert()
    $ws.Rows.Item(10).Insert()
    $ws.Rows.Item(10).Insert()

    $rowsForSheet = $sheetData[$s]

    for ($i = 0; $i -lt 3; $i++) {
        $r = 10 + $i
        $values = $rowsForSheet[$i]

        $ws.Cells.Item($r, 1).Value = $newLabels[$i]
        $ws.Cells.Item($r, 2).Value = $values[0]
        $ws.Cells.Item($r, 3).Value = $values[1]
        $ws.Cells.Item($r, 4).Value = $values[2]
        $ws.Cells.Item($r, 5).Value = $values[3]
        $ws.Cells.Item($r, 6).Value = $values[4]
        $ws.Cells.Item($r, 7).Value = $values[5]
        $ws.Cells.Item($r, 8).Value = $values[6]
        $ws.Cells.Item($r, 9).Value = $values[7]
    }
}

Write-Output "COG Dataset 6 added to all sheets"
